$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price (D) and Volume(1h) (E) columns so that
# values such as "309.49" and "-1.08%" are stored as text, matching the original
# inline-string cell type instead of being auto-converted to numeric/percentage.

$cells = @(
    "D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50","D51","E2","E3","E4","E5","E6","E7","E8","E9","E10","E11","E12","E13","E14","E15","E16","E17","E18","E19","E20","E21","E22","E23","E24","E25","E26","E39","E40","E41","E42","E43","E44","E45","E46","E47","E48","E49","E50","E51"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price column (D)
$ws.Range("D2").Value = "309.49"
$ws.Range("D3").Value = "36.44"
$ws.Range("D4").Value = "5.107"
$ws.Range("D5").Value = "0.07715"
$ws.Range("D6").Value = "8.318"
$ws.Range("D8").Value = "2.986"
$ws.Range("D9").Value = "0.9231"
$ws.Range("D10").Value = "0.1147"
$ws.Range("D11").Value = "0.1880"
$ws.Range("D12").Value = "0.08808"
$ws.Range("D13").Value = "0.03362"
$ws.Range("D14").Value = "0.09534"
$ws.Range("D15").Value = "0.001377"
$ws.Range("D16").Value = "0.005912"
$ws.Range("D17").Value = "3.356"
$ws.Range("D18").Value = "4.390"
$ws.Range("D19").Value = "0.3436"
$ws.Range("D20").Value = "6.318"
$ws.Range("D21").Value = "0.1289"
$ws.Range("D22").Value = "0.2311"
$ws.Range("D23").Value = "0.04336"
$ws.Range("D25").Value = "0.004251"
$ws.Range("D26").Value = "0.0001330"
$ws.Range("D27").Value = "0.0002901"
$ws.Range("D39").Value = "0.02126"
$ws.Range("D40").Value = "0.05006"
$ws.Range("D41").Value = "0.007521"
$ws.Range("D42").Value = "0.1351"
$ws.Range("D43").Value = "0.008391"
$ws.Range("D44").Value = "0.002070"
$ws.Range("D45").Value = "0.007714"
$ws.Range("D46").Value = "0.00006314"
$ws.Range("D48").Value = "0.002867"
$ws.Range("D49").Value = "0.001690"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").Value = "0.0002000"

# Update Volume(1h) column (E)
$ws.Range("E2").Value = "-1.08%"
$ws.Range("E3").Value = "-3.12%"
$ws.Range("E4").Value = "-0.60%"
$ws.Range("E5").Value = "-2.46%"
$ws.Range("E6").Value = "0.60%"
$ws.Range("E7").Value = "-2.79%"
$ws.Range("E8").Value = "4.75%"
$ws.Range("E9").Value = "0.16%"
$ws.Range("E10").Value = "-6.35%"
$ws.Range("E11").Value = "-2.18%"
$ws.Range("E12").Value = "-3.57%"
$ws.Range("E13").Value = "1.42%"
$ws.Range("E14").Value = "-1.10%"
$ws.Range("E15").Value = "-0.56%"
$ws.Range("E16").Value = "2.11%"
$ws.Range("E17").Value = "-4.45%"
$ws.Range("E18").Value = "-0.71%"
$ws.Range("E19").Value = "-0.25%"
$ws.Range("E20").Value = "19.60%"
$ws.Range("E21").Value = "1.34%"
$ws.Range("E22").Value = "-10.88%"
$ws.Range("E23").Value = "-0.72%"
$ws.Range("E24").Value = "-3.98%"
$ws.Range("E25").Value = "-1.24%"
$ws.Range("E26").Value = "8.91%"
$ws.Range("E39").Value = "-4.55%"
$ws.Range("E40").Value = "-2.51%"
$ws.Range("E41").Value = "1.75%"
$ws.Range("E42").Value = "-0.72%"
$ws.Range("E43").Value = "-4.19%"
$ws.Range("E44").Value = "2.88%"
$ws.Range("E45").Value = "-10.80%"
$ws.Range("E46").Value = "-6.37%"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("E48").Value = "-14.75%"
$ws.Range("E49").Value = "40.66%"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("E51").Value = "-0.09%"
